# "cleaned up autograde file"
# Appends the "Next:" follow-up list (and its sub-bullets) to the end of
# the "Progress" list (numId 5), matching the author's new notes, and
# refreshes the header's cached PAGE field now that the document spills
# onto a second page.

$d = $word.ActiveDocument

$newItems = @(
    @{ Text = "Next: "; Level = 1 },
    @{ Text = "error out more gracefully when no answer files"; Level = 2 },
    @{ Text = "more sophisticated grading for errors, i.e. point system, also encoded in a dictionary – can add in a grading method that calls check_hw_answers; will also need more sophisticated records than just whether passed each test wholly"; Level = 2 },
    @{ Text = "need approximate equality checking – may want to make this method more complex, deal with scaling of some kind"; Level = 2 },
    @{ Text = "potential need for more dynamic grader w/ regard to random seeding – see above problem example"; Level = 2 },
    @{ Text = "possible gui? – would only need two buttons + folder navigation; might be nice. "; Level = 2 }
)

foreach ($item in $newItems) {
    $insertAt = $d.Paragraphs.Last.Range
    $insertAt.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $item.Text
    $newPara.Range.ListFormat.ListLevelNumber = $item.Level
}

# The document now overflows onto a second page, so the header's cached
# "PAGE" field result needs to move from 1 to 2.
$header = $d.Sections.Item(1).Headers.Item(1)
$header.Range.Find.Execute("1", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "2", 2)
